$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.011.83'
$ws.Range('E2').Value = '  -3.90%  '
$ws.Range('D3').Value = '1.958.68'
$ws.Range('E3').Value = '  -5.82%  '
$ws.Range('D4').Value = "'1.007"
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').Value = "'326.62"
$ws.Range('E5').Value = '  -3.57%  '
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').Value = "'0.4953"
$ws.Range('E7').Value = '  -5.78%  '
$ws.Range('E8').Value = '  -3.77%  '
$ws.Range('D9').Value = "'52.79"
$ws.Range('E9').Value = '  -3.78%  '
$ws.Range('D10').Value = "'0.09184"
$ws.Range('E10').Value = '  -1.58%  '
$ws.Range('E11').Value = '  -6.41%  '
$ws.Range('D12').Value = "'22.80"
$ws.Range('E12').Value = '  -6.84%  '
$ws.Range('D13').Value = '1.956.41'
$ws.Range('E13').Value = '  -6.71%  '
$ws.Range('D14').Value = "'6.430"
$ws.Range('E14').Value = '  -6.05%  '
$ws.Range('E15').Value = '  -7.51%  '
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').Value = "'91.21"
$ws.Range('E17').Value = '  -9.01%  '
$ws.Range('D18').Value = "'0.00001095"
$ws.Range('E18').Value = '  -5.38%  '
$ws.Range('D19').Value = "'0.06689"
$ws.Range('E19').Value = '  -0.16%  '
$ws.Range('D20').Value = "'19.21"
$ws.Range('E20').Value = '  -7.78%  '
$ws.Range('D21').Value = "'1.005"
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('D22').Value = "'5.942"
$ws.Range('E22').Value = '  -5.57%  '
$ws.Range('D23').Value = '29.055.90'
$ws.Range('E23').Value = '  -3.88%  '
$ws.Range('D24').Value = "'12.00"
$ws.Range('E24').Value = '  -2.97%  '
$ws.Range('D25').Value = "'2.260"
$ws.Range('E25').Value = '  -2.30%  '
$ws.Range('D26').Value = '2.253.07'
$ws.Range('E26').Value = '  -3.23%  '
$ws.Range('D27').Value = "'20.55"
$ws.Range('E27').Value = '  -5.37%  '
$ws.Range('D28').Value = "'155.51"
$ws.Range('E28').Value = '  -4.21%  '
$ws.Range('D29').Value = "'6.258"
$ws.Range('E29').Value = '  -7.71%  '
$ws.Range('D30').Value = "'2.242"
$ws.Range('E30').Value = '  -9.52%  '
$ws.Range('D31').Value = "'126.03"
$ws.Range('E31').Value = '  -5.30%  '
$ws.Range('D32').Value = "'1.039"
$ws.Range('E32').Value = '  -7.78%  '
$ws.Range('D33').Value = "'0.09804"
$ws.Range('E33').Value = '  -6.21%  '
$ws.Range('E34').Value = '  -8.31%  '
$ws.Range('D35').Value = "'5.818"
$ws.Range('E35').Value = '  -6.70%  '
$ws.Range('E36').Value = '  -5.87%  '
$ws.Range('D37').Value = "'0.02417"
$ws.Range('E37').Value = '  -7.07%  '
$ws.Range('D38').Value = "'1.325"
$ws.Range('E38').Value = '  -0.12%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = "'0.06338"
$ws.Range('E39').Value = '  -5.47%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = "'8.969"
$ws.Range('E40').Value = '  -8.85%  '
$ws.Range('D41').Value = "'0.6419"
$ws.Range('E42').Value = '  -9.16%  '
$ws.Range('D43').Value = "'0.1972"
$ws.Range('E43').Value = '  -10.23%  '
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('B45').Value = 'WEMIXTOKEN'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = "'1.363"
$ws.Range('E45').Value = '  +4.21%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = "'0.6199"
$ws.Range('E46').Value = '  -7.55%  '
$ws.Range('D47').Value = "'13.32"
$ws.Range('E47').Value = '  -6.17%  '
$ws.Range('D48').Value = "'2.185"
$ws.Range('E48').Value = '  -7.26%  '
$ws.Range('D49').Value = "'3.458"
$ws.Range('E49').Value = '  -4.58%  '
$ws.Range('E50').Value = '  -5.97%  '
$ws.Range('D51').Value = "'0.07020"
$ws.Range('E51').Value = '  -2.62%  '
